$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The same list of UK equity tickers that already populates the tail of
# column A is appended once more, starting right after the last used row.
$tickers = @("AAF", "ABDN", "ABF", "ANTO", "AUTO", "AV", "BARC", "BATS", "BDEV", "BEZ", "BF.B", "BKG", "BNZL", "BRBY", "BRK.B", "BT-A", "CCH", "CRDA", "DCC", "DGE", "ENT", "EXPN", "FCIT", "FRAS", "GLEN", "HLMA", "HSBA", "HSX", "IMB", "INF", "ITRK", "JMAT", "KGF", "LGEN", "LLOY", "LSEG", "MNDI", "MNG", "OCDO", "PHNX", "PSON", "REL", "RMV", "RR", "RS1", "SBRY", "SDR", "SGRO", "SKG", "SMDS", "SMT", "SN", "SPX", "SSE", "STAN", "STJ", "ULVR", "UU", "WEIR", "WTB")

$startRow = $ws.Cells.Item($ws.UsedRange.Rows.Count, 1).Row + 1

for ($i = 0; $i -lt $tickers.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $tickers[$i]
}
